# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G (header "K") values for rows 2-45 are recalculated; write the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(0,0,1,2,2,0,1,0,0,0,0,0,1,0,2,1,0,1,0,1,0,0,3,0,0,2,1,2,1,2,0,0,2,2,1,0,2,0,1,1,2,2,1,2)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
